$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# The sheet currently has 7 data rows (rows 2-8). Two new students need to be
# inserted: "COBOS" becomes the new row 2, and "TAMAYO" becomes the new row 4
# (right after "MORALES", which shifts from row 2 to row 3).
#
# Rather than using Rows.Insert() (which drags along formatting from the row
# above and leaves stray unused style entries behind), shift the existing
# data down manually by copying values bottom-up, then overwrite the two
# freed-up rows with the new students' data.

# Step 1: shift current rows 2-8 down to rows 3-9, freeing row 2 for COBOS.
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Step 2: shift rows 4-9 down to rows 5-10, freeing row 4 for TAMAYO
# (row 3 now holds MORALES, so the new student is inserted right after it).
for ($r = 9; $r -ge 4; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Row 2: new student COBOS
$ws.Cells.Item(2, 1).Value = 21330051920007
$ws.Cells.Item(2, 2).Value = "COBOS"
$ws.Cells.Item(2, 3).Value = "NOLASCO"
$ws.Cells.Item(2, 4).Value = "YOLET"
$ws.Cells.Item(2, 5).Value = "QUÍMICA I"
$ws.Cells.Item(2, 6).Value = "1AV"
$ws.Cells.Item(2, 7).Value = 6

# Row 4: new student TAMAYO
$ws.Cells.Item(4, 1).Value = 21330051920025
$ws.Cells.Item(4, 2).Value = "TAMAYO"
$ws.Cells.Item(4, 3).Value = "VARGAS"
$ws.Cells.Item(4, 4).Value = "JOSMAR JAHIR"
$ws.Cells.Item(4, 5).Value = "QUÍMICA I"
$ws.Cells.Item(4, 6).Value = "1AV"
$ws.Cells.Item(4, 7).Value = 6
